$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 49, shifting existing rows 49-83 down to 50-84
$ws.Rows.Item(49).Insert()

# Populate the new row 49 with the latest week's data
$ws.Range("A49").Value = 7
$ws.Range("B49").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C49").Value = "Ñuble"
$ws.Range("D49").Value = 44603
$ws.Range("E49").Value = 16
$ws.Range("F49").Value = 100112030
$ws.Range("G49").Value = "Poroto granado"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 60
$ws.Range("K49").Value = 23000
$ws.Range("L49").Value = 24000
$ws.Range("M49").Value = 23500
$ws.Range("N49").Value = '$/saco 25 kilos'
$ws.Range("O49").Value = "Provincia de Diguillín"
$ws.Range("P49").Value = 940
$ws.Range("Q49").Value = 25
$ws.Range("R49").Value = "Hortaliza"
